$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append new status report rows (43-45): date, hours worked, and comment.
$ws.Range("A43").Value = 40260
$ws.Range("B43").Value = 2.5
$ws.Range("C43").Value = "Group Meeting"

$ws.Range("A44").Value = 40260
$ws.Range("B44").Value = 1
$ws.Range("C44").Value = "Weekly Meeting"

$ws.Range("A45").Value = 40263
$ws.Range("B45").Value = 0.5
$ws.Range("C45").Value = "Skype Meeting"

# Match the date formatting used by the rest of the Date column.
$ws.Range("A42").Copy() | Out-Null
$ws.Range("A43:A45").PasteSpecial(-4122) | Out-Null

$ws.Range("A46").Select()
